$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2 through 110: update the date serial
# value from 45172 (2023-09-03) to 45175 (2023-09-06) for every row.
for ($r = 2; $r -le 110; $r++) {
    $ws.Range("C$r").Value = 45175
}
